$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

function Shift-RowLeft {
    param($ws, $row, $newValue)
    $rng = $ws.Range("E$($row):N$($row)")
    $vals = $rng.Value2
    $new = New-Object "object[,]" 1,10
    for ($i = 1; $i -le 9; $i++) {
        $new[0, $i-1] = $vals[1, $i+1]
    }
    $new[0, 9] = $newValue
    $rng.Value2 = $new
}

Shift-RowLeft $ws 8 "فصل چهارم منتهی به 1401/12"
Shift-RowLeft $ws 37 "فصل چهارم منتهی به 1401/12"
Shift-RowLeft $ws 66 "فصل چهارم منتهی به 1401/12"
Shift-RowLeft $ws 88 "فصل چهارم منتهی به 1401/12"
Shift-RowLeft $ws 117 "فصل چهارم منتهی به 1401/12"

Shift-RowLeft $ws 11 69001
Shift-RowLeft $ws 12 555309
Shift-RowLeft $ws 13 "-"
Shift-RowLeft $ws 14 "-"
Shift-RowLeft $ws 15 "-"
Shift-RowLeft $ws 16 53636
Shift-RowLeft $ws 17 "-"
Shift-RowLeft $ws 18 "-"
Shift-RowLeft $ws 19 "-"
Shift-RowLeft $ws 20 "-"
Shift-RowLeft $ws 21 677946
Shift-RowLeft $ws 23 13333
Shift-RowLeft $ws 24 405
Shift-RowLeft $ws 25 153184
Shift-RowLeft $ws 26 "-"
Shift-RowLeft $ws 27 "-"
Shift-RowLeft $ws 28 166922
Shift-RowLeft $ws 30 0
Shift-RowLeft $ws 31 "-"
Shift-RowLeft $ws 32 0
Shift-RowLeft $ws 33 844868
Shift-RowLeft $ws 40 1016535
Shift-RowLeft $ws 41 3797275
Shift-RowLeft $ws 42 "-"
Shift-RowLeft $ws 43 "-"
Shift-RowLeft $ws 44 "-"
Shift-RowLeft $ws 45 171633
Shift-RowLeft $ws 46 "-"
Shift-RowLeft $ws 47 "-"
Shift-RowLeft $ws 48 "-"
Shift-RowLeft $ws 49 "-"
Shift-RowLeft $ws 50 4985443
Shift-RowLeft $ws 52 297292
Shift-RowLeft $ws 53 6529
Shift-RowLeft $ws 54 1401635
Shift-RowLeft $ws 55 "-"
Shift-RowLeft $ws 56 "-"
Shift-RowLeft $ws 57 1705456
Shift-RowLeft $ws 59 0
Shift-RowLeft $ws 60 "-"
Shift-RowLeft $ws 61 0
Shift-RowLeft $ws 62 6690899
Shift-RowLeft $ws 69 14732178
Shift-RowLeft $ws 70 6838130
Shift-RowLeft $ws 71 "-"
Shift-RowLeft $ws 72 "-"
Shift-RowLeft $ws 73 "-"
Shift-RowLeft $ws 74 3199959
Shift-RowLeft $ws 75 "-"
Shift-RowLeft $ws 76 "-"
Shift-RowLeft $ws 77 "-"
Shift-RowLeft $ws 78 "-"
Shift-RowLeft $ws 80 22297457
Shift-RowLeft $ws 81 16120988
Shift-RowLeft $ws 82 9150009
Shift-RowLeft $ws 83 "-"
Shift-RowLeft $ws 84 "-"
Shift-RowLeft $ws 91 -632870
Shift-RowLeft $ws 92 -1838256
Shift-RowLeft $ws 93 "-"
Shift-RowLeft $ws 94 "-"
Shift-RowLeft $ws 95 "-"
Shift-RowLeft $ws 96 -125746
Shift-RowLeft $ws 97 "-"
Shift-RowLeft $ws 98 "-"
Shift-RowLeft $ws 99 "-"
Shift-RowLeft $ws 100 "-"
Shift-RowLeft $ws 101 -2596872
Shift-RowLeft $ws 103 -127881
Shift-RowLeft $ws 104 -2757
Shift-RowLeft $ws 105 -403373
Shift-RowLeft $ws 106 "-"
Shift-RowLeft $ws 107 "-"
Shift-RowLeft $ws 108 -534011
Shift-RowLeft $ws 110 0
Shift-RowLeft $ws 111 "-"
Shift-RowLeft $ws 112 0
Shift-RowLeft $ws 113 -3130883
Shift-RowLeft $ws 120 383665
Shift-RowLeft $ws 121 1959019
Shift-RowLeft $ws 122 "-"
Shift-RowLeft $ws 123 "-"
Shift-RowLeft $ws 124 "-"
Shift-RowLeft $ws 125 45887
Shift-RowLeft $ws 126 "-"
Shift-RowLeft $ws 127 "-"
Shift-RowLeft $ws 128 "-"
Shift-RowLeft $ws 129 "-"
Shift-RowLeft $ws 130 2388571
Shift-RowLeft $ws 132 169411
Shift-RowLeft $ws 133 3772
Shift-RowLeft $ws 134 998262
Shift-RowLeft $ws 135 "-"
Shift-RowLeft $ws 136 "-"
Shift-RowLeft $ws 137 1171445
Shift-RowLeft $ws 139 0
Shift-RowLeft $ws 140 3560016

Write-Host "done"